$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.519122004508972
$ws.Range("B1").Value = 1.816242694854736
$ws.Range("C1").Value = 1.939616203308105
$ws.Range("D1").Value = 1.541638374328613
$ws.Range("E1").Value = 1.350109934806824
